$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 19.35876433333334
$ws.Range("H2").Value = 58.07629300000001
$ws.Range("I2").Value = 0.9707667559429034
$ws.Range("J2").Value = 0.9707667559429034
$ws.Range("M2").Value = 31.618405
$ws.Range("N2").Value = 94.855215
$ws.Range("O2").Value = 0.8578613706944929
$ws.Range("P2").Value = 0.8578613706944929
$ws.Range("Q2").Value = 612.0932509908885
$ws.Range("R2").Value = 5508.839258917996
$ws.Range("S2").Value = 0.8327832998778253
$ws.Range("T2").Value = 0.8327832998778253
$ws.Range("G3").Value = 19.35876433333334
$ws.Range("H3").Value = 58.07629300000001
$ws.Range("I3").Value = 0.9707667559429034
$ws.Range("J3").Value = 0.9707667559429034
$ws.Range("O3").Value = 0.08747555172986397
$ws.Range("P3").Value = 0.08747555172986396
$ws.Range("Q3").Value = 62.41474050429323
$ws.Range("R3").Value = 561.732664538639
$ws.Range("S3").Value = 0.08491835757711567
$ws.Range("T3").Value = 0.08491835757711567
$ws.Range("G4").Value = 19.35876433333334
$ws.Range("H4").Value = 58.07629300000001
$ws.Range("I4").Value = 0.9707667559429034
$ws.Range("J4").Value = 0.9707667559429034
$ws.Range("M4").Value = 2.014730333333334
$ws.Range("N4").Value = 6.044191000000001
$ws.Range("O4").Value = 0.05466307757564324
$ws.Range("P4").Value = 0.05466307757564324
$ws.Range("Q4").Value = 39.00268971821812
$ws.Range("R4").Value = 351.0242074639631
$ws.Range("S4").Value = 0.05306509848796245
$ws.Range("T4").Value = 0.05306509848796245
$ws.Range("H5").Value = 0.9049070000000001
$ws.Range("I5").Value = 0.01512585579145048
$ws.Range("J5").Value = 0.01512585579145048
$ws.Range("M5").Value = 31.618405
$ws.Range("N5").Value = 94.855215
$ws.Range("O5").Value = 0.8578613706944929
$ws.Range("P5").Value = 0.8578613706944929
$ws.Range("Q5").Value = 9.537238671111668
$ws.Range("R5").Value = 85.83514804000501
$ws.Range("S5").Value = 0.01297588738218094
$ws.Range("T5").Value = 0.01297588738218094
$ws.Range("H6").Value = 0.9049070000000001
$ws.Range("I6").Value = 0.01512585579145048
$ws.Range("J6").Value = 0.01512585579145048
$ws.Range("O6").Value = 0.08747555172986397
$ws.Range("P6").Value = 0.08747555172986396
$ws.Range("S6").Value = 0.001323142580743489
$ws.Range("T6").Value = 0.001323142580743489
$ws.Range("H7").Value = 0.9049070000000001
$ws.Range("I7").Value = 0.01512585579145048
$ws.Range("J7").Value = 0.01512585579145048
$ws.Range("M7").Value = 2.014730333333334
$ws.Range("N7").Value = 6.044191000000001
$ws.Range("O7").Value = 0.05466307757564324
$ws.Range("P7").Value = 0.05466307757564324
$ws.Range("Q7").Value = 0.6077145272485557
$ws.Range("R7").Value = 5.469430745237001
$ws.Range("S7").Value = 0.00082682582852605
$ws.Range("T7").Value = 0.0008268258285260501
$ws.Range("G8").Value = 0.2813256666666666
$ws.Range("H8").Value = 0.843977
$ws.Range("I8").Value = 0.01410738826564608
$ws.Range("J8").Value = 0.01410738826564608
$ws.Range("M8").Value = 31.618405
$ws.Range("N8").Value = 94.855215
$ws.Range("O8").Value = 0.8578613706944929
$ws.Range("P8").Value = 0.8578613706944929
$ws.Range("Q8").Value = 8.895068865561665
$ws.Range("R8").Value = 80.055619790055
$ws.Range("S8").Value = 0.01210218343448655
$ws.Range("T8").Value = 0.01210218343448655
$ws.Range("G9").Value = 0.2813256666666666
$ws.Range("H9").Value = 0.843977
$ws.Range("I9").Value = 0.01410738826564608
$ws.Range("J9").Value = 0.01410738826564608
$ws.Range("O9").Value = 0.08747555172986397
$ws.Range("P9").Value = 0.08747555172986396
$ws.Range("Q9").Value = 0.907024238730111
$ws.Range("R9").Value = 8.163218148570998
$ws.Range("S9").Value = 0.0012340515720048
$ws.Range("T9").Value = 0.0012340515720048
$ws.Range("G10").Value = 0.2813256666666666
$ws.Range("H10").Value = 0.843977
$ws.Range("I10").Value = 0.01410738826564608
$ws.Range("J10").Value = 0.01410738826564608
$ws.Range("M10").Value = 2.014730333333334
$ws.Range("N10").Value = 6.044191000000001
$ws.Range("O10").Value = 0.05466307757564324
$ws.Range("P10").Value = 0.05466307757564324
$ws.Range("S10").Value = 0.0007711532591547308
$ws.Range("T10").Value = 0.0007711532591547309
